$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Septiembre de 2020 a las 07:01"

# --- India: updated case counts ---
$ws.Range("B5").Value = 4754356
$ws.Range("C5").Value = 2568
$ws.Range("D5").Value = 3702595
$ws.Range("E5").Value = 973147

# --- Pakistan: updated case counts ---
$ws.Range("B20").Value = 301481
$ws.Range("C20").Value = 526
$ws.Range("D20").Value = 289429
$ws.Range("E20").Value = 5673
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 6379

# --- Uzbekistan: updated case counts ---
$ws.Range("B60").Value = 46850
$ws.Range("C60").Value = 129
$ws.Range("D60").Value = 43511
$ws.Range("E60").Value = 2953
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 386

# --- Tailandia: updated case counts ---
$ws.Range("B128").Value = 3473
$ws.Range("C128").Value = 7
$ws.Range("E128").Value = 103

# --- Birmania (currently row 143): updated case counts. The increase
# pushes it above Reunion and Estonia once the table is re-sorted below. ---
$ws.Range("B143").Value = 2796
$ws.Range("C143").Value = 201
$ws.Range("E143").Value = 2104
$ws.Range("H143").Value = 16

# --- Butan: updated case counts ---
$ws.Range("B187").Value = 244
$ws.Range("C187").Value = 3
$ws.Range("D187").Value = 159
$ws.Range("E187").Value = 85

# Re-sort the data table (rows 4-219) by total cases (column B) descending,
# same as the rest of the sheet, so Birmania's updated total relocates it
# above Reunion/Estonia.
$dataRange = $ws.Range("A4:H219")
$dataRange.Sort($ws.Range("B4"), 2, $null, $null, $null, $null, $null, 0)
